$p = $ppt.ActivePresentation

# --- Slide 8: split "On saute le bloque de code des classes englobantes"
#     into "On saute le " / "bloc " / "de code des classes englobantes"
#     (fixes the typo "bloque" -> "bloc") ---
$s8 = $p.Slides.Item(8)
$shape8 = $s8.Shapes.Item(2)
$para8 = $shape8.TextFrame.TextRange.Paragraphs(3, 1)

# "On saute le bloque de code des classes englobantes"
#  123456789012
$para8.Characters(1, 12).Text = "On saute le "

# remaining original text from char 13: "bloque de code des classes englobantes"
$para8.Characters(13, 7).Text = "bloc "

# --- Slide 9: merge the split runs "ins = C" + "()" into a single run "ins = C()" ---
$s9 = $p.Slides.Item(9)
$shape9 = $s9.Shapes.Item(1)
$para9 = $shape9.TextFrame.TextRange.Paragraphs(9, 1)

$run1 = $para9.Runs(1, 1)
$run1.Text = "ins = C()"
$run2 = $para9.Runs(2, 1)
$run2.Text = ""
